# Apply the changes described by the commit:
# "Add updated spreadsheet and new database table image"
#
# Concretely (content-level changes visible in the sheet):
#   - Column B ("Номер_замовлення" column in the first mini-table) was
#     resized narrower by the user.
#   - Column M ("id_замовлення" header of the "3 НФ" table) was widened
#     with a best-fit/auto-fit so its header text is no longer truncated.
#   - The "3 НФ" table area (selection) became the active selection, with
#     the view scrolled down to show rows/cols around it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure this is the active/selected sheet in the workbook.
$ws.Activate()

# --- Column B: manual narrow resize ---
$ws.Columns.Item(2).ColumnWidth = 5.6667

# --- Column M (13): best-fit / auto-fit to its (now visible) header text ---
$ws.Columns.Item(13).AutoFit()

# --- View state: scroll so row 16 / column E is the top-left corner, and
#     select the range below the "3 НФ" table (L31:Q31), matching where the
#     user ended up after adding the new table image. ---
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("L31:Q31").Select()
